$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Efna1","Epha4","ECs",3,1,20.447252,61.341756,0.8699145605694745,0.8770588936480435,3,1,7.727270999999999,23.181813,0.4492515923977304,0.4784711627054499,158.001457409292,1422.013116683628,0.3908105015858083,0.4196473886049349),
    @("ECs","Efna1","Epha4","FAPs",3,1,20.447252,61.341756,0.8699145605694745,0.8770588936480435,3,1,6.292848333333333,18.878545,0.3658564756519351,0.3896519817642027,128.6714556694467,1158.04310102502,0.3182638752482497,0.3417477360338793),
    @("ECs","Efna1","Epha4","M1",3,1,20.447252,61.341756,0.8699145605694745,0.8770588936480435,1,0.3333333333333333,0.021228,0.063684,0.001234163109255392,0.001314433755709007,0.434054265456,3.906488389104,0.001073616458858961,0.001152835815555784),
    @("ECs","Efna1","Epha4","M2",3,1,20.447252,61.341756,0.8699145605694745,0.8770588936480435,1,0.3333333333333333,0.007773333333333333,0.02332,0.0004519295852621654,0.0004813233336965963,0.1589433055466667,1.43048974992,0.0003931401265716815,0.0004221489105389248),
    @("ECs","Efna1","Epha4","sCs",3,1,20.447252,61.341756,0.8699145605694745,0.8770588936480435,2,1,3.151199,6.302398,0.183205839255817,0.1300810984409417,64.433360055148,386.600160330888,0.1593734271499858,0.1140887842831346),
    @("FAPs","Efna1","Epha4","ECs",3,1,2.483247333333333,7.449742,0.105648084777455,0.1065157390747562,3,1,7.727270999999999,23.181813,0.4492515923977304,0.4784711627054499,19.188725104694,172.698525942246,0.04746257032004208,0.05096470952152894),
    @("FAPs","Efna1","Epha4","FAPs",3,1,2.483247333333333,7.449742,0.105648084777455,0.1065157390747562,3,1,6.292848333333333,18.878545,0.3658564756519351,0.3896519817642027,15.62669884282111,140.64028958539,0.03865203595605653,0.0415040688195575),
    @("FAPs","Efna1","Epha4","M1",3,1,2.483247333333333,7.449742,0.105648084777455,0.1065157390747562,1,0.3333333333333333,0.021228,0.063684,0.001234163109255392,0.001314433755709007,0.052714374392,0.474429369528,0.0001303869687958211,0.0001400078829541524),
    @("FAPs","Efna1","Epha4","M2",3,1,2.483247333333333,7.449742,0.105648084777455,0.1065157390747562,1,0.3333333333333333,0.007773333333333333,0.02332,0.0004519295852621654,0.0004813233336965963,0.01930310927111111,0.17372798344,0.00004774549513721733,0.00005126851062261847),
    @("FAPs","Efna1","Epha4","sCs",3,1,2.483247333333333,7.449742,0.105648084777455,0.1065157390747562,2,1,3.151199,6.302398,0.183205839255817,0.1300810984409417,7.825206513552667,46.951239081316,0.01935534603742334,0.01385568434009303),
    @("sCs","Efna1","Epha4","ECs",2,1,0.5743975,1.148795,0.02443735465307048,0.01642536727720028,3,1,7.727270999999999,23.181813,0.4492515923977304,0.4784711627054499,4.438525144222499,26.631150865335,0.01097852049188,0.00785906457898607),
    @("sCs","Efna1","Epha4","FAPs",2,1,0.5743975,1.148795,0.02443735465307048,0.01642536727720028,3,1,6.292848333333333,18.878545,0.3658564756519351,0.3896519817642027,3.614596350545833,21.687578103275,0.008940564447628785,0.006400176910765977),
    @("sCs","Efna1","Epha4","M1",2,1,0.5743975,1.148795,0.02443735465307048,0.01642536727720028,1,0.3333333333333333,0.021228,0.063684,0.001234163109255392,0.001314433755709007,0.01219331013,0.07315986078,0.00003015968160061019,0.00002159005719907019),
    @("sCs","Efna1","Epha4","M2",2,1,0.5743975,1.148795,0.02443735465307048,0.01642536727720028,1,0.3333333333333333,0.007773333333333333,0.02332,0.0004519295852621654,0.0004813233336965963,0.004464983233333333,0.0267898994,0.00001104396355326659,0.000007905912535053025),
    @("sCs","Efna1","Epha4","sCs",2,1,0.5743975,1.148795,0.02443735465307048,0.01642536727720028,2,1,3.151199,6.302398,0.183205839255817,0.1300810984409417,1.8100408276025,7.24016331041,0.004477066068407822,0.002136629817714113)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

Write-Output "Updated rows 2..16"
